# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets,
# matching the "2026-02-06" 09:59-10:00 Bathroom readings recorded by the
# logger. Values are forced to Text (NumberFormat "@") on the Date column
# (and the Value column, which sometimes holds "NN.N%" strings) before the
# write so Excel's autodetect does not silently turn them into date
# serials / percentage numbers.

$wb = $excel.ActiveWorkbook

function Append-LogRows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        $Rows
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $r = $StartRow
    foreach ($row in $Rows) {
        $ws.Cells.Item($r, 1).NumberFormat = "@"
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).NumberFormat = "@"
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $r++
    }
}

# PIR sheet: rows 231-243 (dimension A1:F230 -> A1:F243)
$pirRows = @(
    @("2026-02-06","09:59:16","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:59:16","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:59:21","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:59:26","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:59:31","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:59:36","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:59:41","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:59:46","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:59:51","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:59:56","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:00:01","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:00:06","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:00:11","10:00","Bathroom","No Motion","Inactive")
)
Append-LogRows "PIR" 231 $pirRows

# Humidity sheet: rows 142-150 (dimension A1:F141 -> A1:F150)
$humidityRows = @(
    @("2026-02-06","09:59:18","09:00","Bathroom","69.8%","Active"),
    @("2026-02-06","09:59:29","09:00","Bathroom","70.0%","Active"),
    @("2026-02-06","09:59:39","09:00","Bathroom","69.6%","Active"),
    @("2026-02-06","09:59:44","09:00","Bathroom","69.8%","Active"),
    @("2026-02-06","09:59:49","09:00","Bathroom","69.6%","Active"),
    @("2026-02-06","09:59:59","09:00","Bathroom","69.7%","Active"),
    @("2026-02-06","10:00:04","10:00","Bathroom","69.4%","Active"),
    @("2026-02-06","10:00:09","10:00","Bathroom","69.5%","Active"),
    @("2026-02-06","10:00:14","10:00","Bathroom","69.4%","Active")
)
Append-LogRows "Humidity" 142 $humidityRows

# Temperature sheet: rows 142-150 (dimension A1:F141 -> A1:F150)
$temperatureRows = @(
    @("2026-02-06","09:59:19","09:00","Bathroom","27.8C","Active"),
    @("2026-02-06","09:59:29","09:00","Bathroom","27.9C","Active"),
    @("2026-02-06","09:59:39","09:00","Bathroom","27.7C","Active"),
    @("2026-02-06","09:59:44","09:00","Bathroom","27.9C","Active"),
    @("2026-02-06","09:59:49","09:00","Bathroom","27.8C","Active"),
    @("2026-02-06","09:59:59","09:00","Bathroom","27.9C","Active"),
    @("2026-02-06","10:00:04","10:00","Bathroom","27.8C","Active"),
    @("2026-02-06","10:00:09","10:00","Bathroom","27.9C","Active"),
    @("2026-02-06","10:00:15","10:00","Bathroom","27.8C","Active")
)
Append-LogRows "Temperature" 142 $temperatureRows

Write-Output "Appended $($pirRows.Count) PIR rows, $($humidityRows.Count) Humidity rows, $($temperatureRows.Count) Temperature rows."
